$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 288, pushing the existing row 288 (and everything below)
# down by one. This is how a new weekly price record is added to the top of
# the (date-descending) series for this market/product.
$ws.Rows("288:288").Insert()

# Populate the newly inserted row 288 with the new record's data.
$ws.Range("A288").Value = 5
$ws.Range("B288").Value = "Macroferia Regional de Talca"
$ws.Range("C288").Value = "Maule"
$ws.Range("D288").Value = 45027
$ws.Range("E288").Value = 7
$ws.Range("F288").Value = 100112009
$ws.Range("G288").Value = "Acelga"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 500
$ws.Range("K288").Value = 2500
$ws.Range("L288").Value = 2500
$ws.Range("M288").Value = 2500
$ws.Range("N288").Value = "$/docena de atados (4 kilos)"
$ws.Range("O288").Value = "Región del Maule"
$ws.Range("P288").Value = 625
$ws.Range("Q288").Value = 4
$ws.Range("R288").Value = "Hortaliza"
